# Hortaliza, Feria Lagunitas de Puerto Montt - Ají
# Insert two new weekly price observations at row 150, pushing the
# existing row 150 (and everything below it, through row 186) down by
# two rows (old row 150 -> 152, old row 151 -> 153, ..., old row 186 ->
# 188). The two freshly inserted rows (150 and 151) are then populated
# with the new records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows starting at row 150; this shifts old rows 150-186
# down to 152-188, leaving two fresh blank rows at 150 and 151 to fill in.
$ws.Rows.Item(150).Resize(2).EntireRow.Insert()

# --- New row 150 -------------------------------------------------------
$ws.Range("A150").Value = 4
$ws.Range("B150").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C150").Value = "Los Lagos"
$ws.Range("D150").Value = 44551
$ws.Range("E150").Value = 10
$ws.Range("F150").Value = 100112021
$ws.Range("G150").Value = "Ají"
$ws.Range("H150").Value = "Inferno"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 50
$ws.Range("K150").Value = 20000
$ws.Range("L150").Value = 20000
$ws.Range("M150").Value = 20000
$ws.Range("N150").Value = "$/caja 12 kilos"
$ws.Range("O150").Value = "Región de Arica y Parinacota"
$ws.Range("P150").Value = 1667
$ws.Range("Q150").Value = 12
$ws.Range("R150").Value = "Hortaliza"

# --- New row 151 -------------------------------------------------------
$ws.Range("A151").Value = 4
$ws.Range("B151").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C151").Value = "Los Lagos"
$ws.Range("D151").Value = 44551
$ws.Range("E151").Value = 10
$ws.Range("F151").Value = 100112021
$ws.Range("G151").Value = "Ají"
$ws.Range("H151").Value = "Inferno"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 80
$ws.Range("K151").Value = 30000
$ws.Range("L151").Value = 30000
$ws.Range("M151").Value = 30000
$ws.Range("N151").Value = "$/caja 15 kilos"
$ws.Range("O151").Value = "Provincia de Huasco"
$ws.Range("P151").Value = 2000
$ws.Range("Q151").Value = 15
$ws.Range("R151").Value = "Hortaliza"

# Keep the date columns formatted the same way as the rest of column D.
$ws.Range("D150:D151").NumberFormat = $ws.Range("D149").NumberFormat
